# Generate Report for Handoff
# Updates the localization-status report: rows 4-7 (the "Ready for handoff"
# files) get a fresh handoff pass -- Priority flips from "low" to "ht" and
# the handoff timestamps advance a few seconds on the Overview + per-locale
# sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" (col G) for rows 4-7
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in 4..7) {
    $overview.Cells.Item($r, 7).Value = "2016-08-24 18:39:37"
}

# zh-cn sheet: Priority (col E) -> "ht", Latest Handoff Datetime (col H)
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in 4..7) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-24 18:39:32"
}

# de-de sheet: Priority (col E) -> "ht", Latest Handoff Datetime (col H)
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in 4..7) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-24 18:39:37"
}
